# Apply "Penalty Reward System" forecast refresh edit.
# The forecast window rolled forward by one week: each week's start date
# advances by 7 days, and the MyForecast (column D) values are recomputed
# for most weeks. The Summary sheet's derived metrics are updated to match.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Helper: write a value to a cell while forcing text storage (the source
# workbook keeps these columns as plain text, even when the text looks like
# a date or a number).
function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# --- Sheet "Forecast Comparison": Week_Start_Date (B) and MyForecast (D) ---

Set-TextCell $wsForecast "B2" "2025-01-12"
$wsForecast.Range("D2").Value = 31

Set-TextCell $wsForecast "B3" "2025-01-19"
$wsForecast.Range("D3").Value = 34

Set-TextCell $wsForecast "B4" "2025-01-26"
$wsForecast.Range("D4").Value = 35

Set-TextCell $wsForecast "B5" "2025-02-02"
$wsForecast.Range("D5").Value = 36

Set-TextCell $wsForecast "B6" "2025-02-09"
$wsForecast.Range("D6").Value = 26

Set-TextCell $wsForecast "B7" "2025-02-16"
$wsForecast.Range("D7").Value = 28

Set-TextCell $wsForecast "B8" "2025-02-23"
# D8 unchanged (30)

Set-TextCell $wsForecast "B9" "2025-03-02"
$wsForecast.Range("D9").Value = 31

Set-TextCell $wsForecast "B10" "2025-03-09"
$wsForecast.Range("D10").Value = 29

Set-TextCell $wsForecast "B11" "2025-03-16"
$wsForecast.Range("D11").Value = 30

Set-TextCell $wsForecast "B12" "2025-03-23"
$wsForecast.Range("D12").Value = 29

Set-TextCell $wsForecast "B13" "2025-03-30"
$wsForecast.Range("D13").Value = 31

Set-TextCell $wsForecast "B14" "2025-04-06"
$wsForecast.Range("D14").Value = 30

Set-TextCell $wsForecast "B15" "2025-04-13"
$wsForecast.Range("D15").Value = 29

Set-TextCell $wsForecast "B16" "2025-04-20"
$wsForecast.Range("D16").Value = 30

Set-TextCell $wsForecast "B17" "2025-04-27"
$wsForecast.Range("D17").Value = 29

# --- Sheet "Summary": refreshed derived metrics ---

Set-TextCell $wsSummary "B2"  "2023-02-12 to 2025-01-05"
Set-TextCell $wsSummary "B4"  "84"
Set-TextCell $wsSummary "B8"  "2713 units"
Set-TextCell $wsSummary "B9"  "486"
Set-TextCell $wsSummary "B10" "251"
Set-TextCell $wsSummary "B11" "136"
Set-TextCell $wsSummary "B12" "36"
Set-TextCell $wsSummary "B13" "2025-02-02"
Set-TextCell $wsSummary "B14" "26"
